# 414: Add ARMS tab to test extract files
$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (CMS) so it lands at the end,
# and becomes the active/selected tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "ARMS"

$headers = @(
    "Assessment_Staff_Name",
    "Assessment_Staff_Key",
    "Assessment_Staff_Grade",
    "Assessmentent_Team_Key",
    "Assessment_Provider_Code",
    "CRN",
    "Disposal_or_Release_Date",
    "Sentence Type",
    "SO_Registration_Date"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$headerRange = $newSheet.Range("A1:I1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 9
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 0xFFFFFF
$headerRange.Interior.Color = 0xA0640B
$headerRange.Interior.PatternColor = 0xFFFFFF
$headerRange.HorizontalAlignment = -4131
$headerRange.WrapText = $true
$headerRange.RowHeight = 37

# Column C (Assessment_Staff_Grade) is formatted as text and does not wrap.
$gradeCell = $newSheet.Range("C1")
$gradeCell.WrapText = $false
$gradeCell.NumberFormat = "@"

# Borders: thin blue lines (top/right) with a slightly muted bottom edge, applied
# per column-group to mirror the header banner look.
$blue = 0xA67738
$mutedBottom = 0xB1A5A5

$colA = $newSheet.Range("A1")
$colA.Borders.Item(7).LineStyle = 1
$colA.Borders.Item(7).Color = $blue
$colA.Borders.Item(10).LineStyle = 1
$colA.Borders.Item(10).Color = $blue
$colA.Borders.Item(8).LineStyle = 1
$colA.Borders.Item(8).Color = $blue
$colA.Borders.Item(9).LineStyle = 1
$colA.Borders.Item(9).Color = $mutedBottom

$colB = $newSheet.Range("B1")
$colB.Borders.Item(10).LineStyle = 1
$colB.Borders.Item(10).Color = $blue
$colB.Borders.Item(8).LineStyle = 1
$colB.Borders.Item(8).Color = $blue
$colB.Borders.Item(9).LineStyle = 1
$colB.Borders.Item(9).Color = $mutedBottom

$colD = $newSheet.Range("D1")
$colD.Borders.Item(10).LineStyle = 1
$colD.Borders.Item(10).Color = $blue
$colD.Borders.Item(8).LineStyle = 1
$colD.Borders.Item(8).Color = $blue
$colD.Borders.Item(9).LineStyle = 1
$colD.Borders.Item(9).Color = $mutedBottom

$colE = $newSheet.Range("E1")
$colE.Borders.Item(10).LineStyle = 1
$colE.Borders.Item(10).Color = $blue
$colE.Borders.Item(8).LineStyle = 1
$colE.Borders.Item(8).Color = $blue
$colE.Borders.Item(9).LineStyle = 1
$colE.Borders.Item(9).Color = $mutedBottom

$colC = $newSheet.Range("C1")
$colC.Borders.Item(10).LineStyle = 1
$colC.Borders.Item(10).Color = $blue
$colC.Borders.Item(8).LineStyle = 1
$colC.Borders.Item(8).Color = $blue
$colC.Borders.Item(9).LineStyle = 1
$colC.Borders.Item(9).Color = $mutedBottom

foreach ($colLetter in @("F", "G", "H", "I")) {
    $cell = $newSheet.Range($colLetter + "1")
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(10).Color = $blue
}

$newSheet.Range("A1:I1").Select()
